# This script updates the NATMI ligand-receptor (Tgfb1-Tgfbr2) expression/
# specificity table with newly recomputed TPM-derived values.
#
# Columns G/H (ligand average/total expression) and I/J (their specificity
# scores) depend only on the "Sending cluster" (col A); columns M/N (receptor
# average/total expression) and O/P (their specificity scores) depend only on
# the "Target cluster" (col D). Columns Q/R (edge expression weights) and S/T
# (edge specificity scores) are derived per row from the above.
#
# All values below are the recalculated results from the updated TPM input;
# they are written directly into the corresponding data cells (rows 2-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 54.53585066666667
$ws.Range("H2").Value = 163.607552
$ws.Range("I2").Value = 0.3031388658437607
$ws.Range("J2").Value = 0.3031388658437607
$ws.Range("M2").Value = 41.83357366666667
$ws.Range("N2").Value = 125.500721
$ws.Range("O2").Value = 0.2773195847425811
$ws.Range("P2").Value = 0.2773195847425811
$ws.Range("Q2").Value = 2281.429526338332
$ws.Range("R2").Value = 20532.86573704499
$ws.Range("S2").Value = 0.08406634439512871
$ws.Range("T2").Value = 0.08406634439512871

# Row 3
$ws.Range("G3").Value = 54.53585066666667
$ws.Range("H3").Value = 163.607552
$ws.Range("I3").Value = 0.3031388658437607
$ws.Range("J3").Value = 0.3031388658437607
$ws.Range("O3").Value = 0.4239803668761465
$ws.Range("P3").Value = 0.4239803668761465
$ws.Range("Q3").Value = 3487.966161772771
$ws.Range("R3").Value = 31391.69545595494
$ws.Range("S3").Value = 0.1285249275548566
$ws.Range("T3").Value = 0.1285249275548566

# Row 4
$ws.Range("G4").Value = 54.53585066666667
$ws.Range("H4").Value = 163.607552
$ws.Range("I4").Value = 0.3031388658437607
$ws.Range("J4").Value = 0.3031388658437607
$ws.Range("M4").Value = 14.18032733333333
$ws.Range("N4").Value = 42.540982
$ws.Range("O4").Value = 0.09400302539123752
$ws.Range("P4").Value = 0.09400302539123752
$ws.Range("Q4").Value = 773.3362138551182
$ws.Range("R4").Value = 6960.025924696064
$ws.Range("S4").Value = 0.02849597050298198
$ws.Range("T4").Value = 0.02849597050298198

# Row 5
$ws.Range("G5").Value = 54.53585066666667
$ws.Range("H5").Value = 163.607552
$ws.Range("I5").Value = 0.3031388658437607
$ws.Range("J5").Value = 0.3031388658437607
$ws.Range("M5").Value = 30.87848266666667
$ws.Range("N5").Value = 92.635448
$ws.Range("O5").Value = 0.204697022990035
$ws.Range("P5").Value = 0.2046970229900349
$ws.Range("Q5").Value = 1683.984319522589
$ws.Range("R5").Value = 15155.8588757033
$ws.Range("S5").Value = 0.0620516233907934
$ws.Range("T5").Value = 0.0620516233907934

# Row 6
$ws.Range("I6").Value = 0.1026363515063155
$ws.Range("J6").Value = 0.1026363515063155
$ws.Range("M6").Value = 41.83357366666667
$ws.Range("N6").Value = 125.500721
$ws.Range("O6").Value = 0.2773195847425811
$ws.Range("P6").Value = 0.2773195847425811
$ws.Range("Q6").Value = 772.4433557880828
$ws.Range("R6").Value = 6951.990202092746
$ws.Range("S6").Value = 0.028463070379225
$ws.Range("T6").Value = 0.028463070379225

# Row 7
$ws.Range("I7").Value = 0.1026363515063155
$ws.Range("J7").Value = 0.1026363515063155
$ws.Range("O7").Value = 0.4239803668761465
$ws.Range("P7").Value = 0.4239803668761465
$ws.Range("S7").Value = 0.04351579796647677
$ws.Range("T7").Value = 0.04351579796647678

# Row 8
$ws.Range("I8").Value = 0.1026363515063155
$ws.Range("J8").Value = 0.1026363515063155
$ws.Range("M8").Value = 14.18032733333333
$ws.Range("N8").Value = 42.540982
$ws.Range("O8").Value = 0.09400302539123752
$ws.Range("P8").Value = 0.09400302539123752
$ws.Range("Q8").Value = 261.8351403303924
$ws.Range("R8").Value = 2356.516262973532
$ws.Range("S8").Value = 0.009648127556712155
$ws.Range("T8").Value = 0.009648127556712155

# Row 9
$ws.Range("I9").Value = 0.1026363515063155
$ws.Range("J9").Value = 0.1026363515063155
$ws.Range("M9").Value = 30.87848266666667
$ws.Range("N9").Value = 92.635448
$ws.Range("O9").Value = 0.204697022990035
$ws.Range("P9").Value = 0.2046970229900349
$ws.Range("Q9").Value = 570.1611572259609
$ws.Range("R9").Value = 5131.450415033648
$ws.Range("S9").Value = 0.02100935560390157
$ws.Range("T9").Value = 0.02100935560390157

# Row 10
$ws.Range("G10").Value = 12.55635966666667
$ws.Range("H10").Value = 37.669079
$ws.Range("I10").Value = 0.06979483370938171
$ws.Range("J10").Value = 0.06979483370938172
$ws.Range("M10").Value = 41.83357366666667
$ws.Range("N10").Value = 125.500721
$ws.Range("O10").Value = 0.2773195847425811
$ws.Range("P10").Value = 0.2773195847425811
$ws.Range("Q10").Value = 525.277397100662
$ws.Range("R10").Value = 4727.496573905959
$ws.Range("S10").Value = 0.01935547430146324
$ws.Range("T10").Value = 0.01935547430146324

# Row 11
$ws.Range("G11").Value = 12.55635966666667
$ws.Range("H11").Value = 37.669079
$ws.Range("I11").Value = 0.06979483370938171
$ws.Range("J11").Value = 0.06979483370938172
$ws.Range("O11").Value = 0.4239803668761465
$ws.Range("P11").Value = 0.4239803668761465
$ws.Range("Q11").Value = 803.0709541888707
$ws.Range("R11").Value = 7227.638587699837
$ws.Range("S11").Value = 0.02959163920216329
$ws.Range("T11").Value = 0.0295916392021633

# Row 12
$ws.Range("G12").Value = 12.55635966666667
$ws.Range("H12").Value = 37.669079
$ws.Range("I12").Value = 0.06979483370938171
$ws.Range("J12").Value = 0.06979483370938172
$ws.Range("M12").Value = 14.18032733333333
$ws.Range("N12").Value = 42.540982
$ws.Range("O12").Value = 0.09400302539123752
$ws.Range("P12").Value = 0.09400302539123752
$ws.Range("Q12").Value = 178.0532901883975
$ws.Range("R12").Value = 1602.479611695578
$ws.Range("S12").Value = 0.006560925525360209
$ws.Range("T12").Value = 0.00656092552536021

# Row 13
$ws.Range("G13").Value = 12.55635966666667
$ws.Range("H13").Value = 37.669079
$ws.Range("I13").Value = 0.06979483370938171
$ws.Range("J13").Value = 0.06979483370938172
$ws.Range("M13").Value = 30.87848266666667
$ws.Range("N13").Value = 92.635448
$ws.Range("O13").Value = 0.204697022990035
$ws.Range("P13").Value = 0.2046970229900349
$ws.Range("Q13").Value = 387.7213343235991
$ws.Range("R13").Value = 3489.492008912392
$ws.Range("S13").Value = 0.01428679468039497
$ws.Range("T13").Value = 0.01428679468039498

# Row 14
$ws.Range("G14").Value = 94.34696966666667
$ws.Range("H14").Value = 283.040909
$ws.Range("I14").Value = 0.524429948940542
$ws.Range("J14").Value = 0.5244299489405421
$ws.Range("M14").Value = 41.83357366666667
$ws.Range("N14").Value = 125.500721
$ws.Range("O14").Value = 0.2773195847425811
$ws.Range("P14").Value = 0.2773195847425811
$ws.Range("Q14").Value = 3946.870905777265
$ws.Range("R14").Value = 35521.83815199539
$ws.Range("S14").Value = 0.1454346956667641
$ws.Range("T14").Value = 0.1454346956667641

# Row 15
$ws.Range("G15").Value = 94.34696966666667
$ws.Range("H15").Value = 283.040909
$ws.Range("I15").Value = 0.524429948940542
$ws.Range("J15").Value = 0.5244299489405421
$ws.Range("O15").Value = 0.4239803668761465
$ws.Range("P15").Value = 0.4239803668761465
$ws.Range("Q15").Value = 6034.178134939677
$ws.Range("R15").Value = 54307.60321445709
$ws.Range("S15").Value = 0.2223480021526497
$ws.Range("T15").Value = 0.2223480021526498

# Row 16
$ws.Range("G16").Value = 94.34696966666667
$ws.Range("H16").Value = 283.040909
$ws.Range("I16").Value = 0.524429948940542
$ws.Range("J16").Value = 0.5244299489405421
$ws.Range("M16").Value = 14.18032733333333
$ws.Range("N16").Value = 42.540982
$ws.Range("O16").Value = 0.09400302539123752
$ws.Range("P16").Value = 0.09400302539123752
$ws.Range("Q16").Value = 1337.870912781404
$ws.Range("R16").Value = 12040.83821503264
$ws.Range("S16").Value = 0.04929800180618316
$ws.Range("T16").Value = 0.04929800180618317

# Row 17
$ws.Range("G17").Value = 94.34696966666667
$ws.Range("H17").Value = 283.040909
$ws.Range("I17").Value = 0.524429948940542
$ws.Range("J17").Value = 0.5244299489405421
$ws.Range("M17").Value = 30.87848266666667
$ws.Range("N17").Value = 92.635448
$ws.Range("O17").Value = 0.204697022990035
$ws.Range("P17").Value = 0.2046970229900349
$ws.Range("Q17").Value = 2913.291267504692
$ws.Range("R17").Value = 26219.62140754223
$ws.Range("S17").Value = 0.107349249314945
$ws.Range("T17").Value = 0.107349249314945
